# Trange instead of T50
# Update p-values in the "Cod" column of the Fig3 p-values table.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "0.72" "0.76"
Replace-Text "0.39" "0.34"
Replace-Text "0.93" "0.97"
Replace-Text "0.04" "0.07"
